$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '20.571.79'
$ws.Range('E2').Value = '  +2.31%  '
$ws.Range('D3').Value = '1.472.65'
$ws.Range('E3').Value = '  +2.99%  '
$ws.Range('D4').Value = '1.008'
$ws.Range('E4').Value = '  +0.64%  '
$ws.Range('D5').Value = '0.9910'
$ws.Range('E5').Value = '  -1.10%  '
$ws.Range('D6').Value = '280.77'
$ws.Range('E6').Value = '  +2.05%  '
$ws.Range('D7').Value = '0.3738'
$ws.Range('E7').Value = '  +1.16%  '
$ws.Range('D8').Value = '0.3216'
$ws.Range('E8').Value = '  +4.59%  '
$ws.Range('D9').Value = '41.72'
$ws.Range('E9').Value = '  +4.00%  '
$ws.Range('D10').Value = '1.073'
$ws.Range('E10').Value = '  +6.87%  '
$ws.Range('D11').Value = '0.06749'
$ws.Range('E11').Value = '  +2.71%  '
$ws.Range('D12').Value = '1.001'
$ws.Range('E12').Value = '  -0.11%  '
$ws.Range('D13').Value = '5.679'
$ws.Range('E13').Value = '  +4.61%  '
$ws.Range('D14').Value = '18.59'
$ws.Range('E14').Value = '  +7.82%  '
$ws.Range('D15').Value = '6.351'
$ws.Range('E15').Value = '  +2.75%  '
$ws.Range('D16').Value = '1.472.27'
$ws.Range('E16').Value = '  +2.82%  '
$ws.Range('D17').Value = '0.00001044'
$ws.Range('E17').Value = '  +3.01%  '
$ws.Range('D18').Value = '0.05823'
$ws.Range('E18').Value = '  -0.07%  '
$ws.Range('D19').Value = '73.36'
$ws.Range('E19').Value = '  -2.73%  '
$ws.Range('D20').Value = '0.9909'
$ws.Range('E20').Value = '  -1.08%  '
$ws.Range('D21').Value = '5.752'
$ws.Range('E21').Value = '  +1.35%  '
$ws.Range('D22').Value = '15.04'
$ws.Range('E22').Value = '  +3.68%  '
$ws.Range('D23').Value = '11.29'
$ws.Range('E23').Value = '  +1.70%  '
$ws.Range('D24').Value = '2.312'
$ws.Range('E24').Value = '  -0.31%  '
$ws.Range('D25').Value = '20.649.61'
$ws.Range('E25').Value = '  +2.65%  '
$ws.Range('D26').Value = '2.351'
$ws.Range('E26').Value = '  +2.41%  '
$ws.Range('D27').Value = '138.91'
$ws.Range('E27').Value = '  +0.33%  '
$ws.Range('D28').Value = '17.79'
$ws.Range('E28').Value = '  +5.22%  '
$ws.Range('D29').Value = '1.641.93'
$ws.Range('E29').Value = '  +3.20%  '
$ws.Range('D30').Value = '114.77'
$ws.Range('E30').Value = '  +4.78%  '
$ws.Range('D31').Value = '3.980'
$ws.Range('E31').Value = '  +1.48%  '
$ws.Range('D32').Value = '5.426'
$ws.Range('E32').Value = '  -0.38%  '
$ws.Range('D33').Value = '0.8560'
$ws.Range('E33').Value = '  -6.78%  '
$ws.Range('D34').Value = '0.07911'
$ws.Range('E34').Value = '  +1.61%  '
$ws.Range('D35').Value = '1.596'
$ws.Range('E35').Value = '  +23.75%  '
$ws.Range('D36').Value = '0.06055'
$ws.Range('E36').Value = '  +6.38%  '
$ws.Range('D37').Value = '4.987'
$ws.Range('E37').Value = '  +3.93%  '
$ws.Range('D38').Value = '10.85'
$ws.Range('E38').Value = '  -4.75%  '
$ws.Range('B39').Value = 'Frax'
$ws.Range('C39').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D39').Value = '0.9934'
$ws.Range('E39').Value = '  -0.79%  '
$ws.Range('D40').Value = '0.02093'
$ws.Range('E40').Value = '  +4.14%  '
$ws.Range('B41').Value = 'FraxShare'
$ws.Range('C41').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D41').Value = '7.852'
$ws.Range('E41').Value = '  -6.88%  '
$ws.Range('B42').Value = 'TrustWalletToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D42').Value = '1.139'
$ws.Range('E42').Value = '  +2.21%  '
$ws.Range('B43').Value = 'Algorand'
$ws.Range('C43').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D43').Value = '0.1930'
$ws.Range('E43').Value = '  +0.41%  '
$ws.Range('D44').Value = '0.5487'
$ws.Range('E44').Value = '  +3.08%  '
$ws.Range('D45').Value = '12.64'
$ws.Range('E45').Value = '  +3.35%  '
$ws.Range('D46').Value = '3.606'
$ws.Range('E46').Value = '  +1.38%  '
$ws.Range('D47').Value = '0.5430'
$ws.Range('E47').Value = '  +5.79%  '
$ws.Range('D48').Value = '121.31'
$ws.Range('E48').Value = '  +9.81%  '
$ws.Range('D49').Value = '1.851'
$ws.Range('E49').Value = '  +3.68%  '
$ws.Range('D50').Value = '1.066'
$ws.Range('E50').Value = '  +1.27%  '
$ws.Range('D51').Value = '0.06461'
$ws.Range('E51').Value = '  +4.00%  '
